$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.296.16"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +0.75%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'1.680.69"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +0.88%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  +0.11%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'218.39"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +0.80%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'0.5254"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  +3.10%  "
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = "'  +0.08%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.2702"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  +2.72%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.06482"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  +1.39%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'21.99"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  +2.13%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.07533"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  +1.83%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'4.528"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  +0.48%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'1.675.13"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +0.02%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'0.5805"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +0.26%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'0.000008544"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  +0.14%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'64.73"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +1.05%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'26.331.31"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +0.65%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'4.924"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +0.06%  "
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = "'  +0.11%  "
$ws.Range("E19").ClearFormats()
$ws.Range("E20").Value = "'  +0.69%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'189.99"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  +0.31%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'6.205"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +0.23%  "
$ws.Range("E22").ClearFormats()
$ws.Range("E23").Value = "'  +0.07%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'145.67"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  +0.35%  "
$ws.Range("E24").ClearFormats()
$ws.Range("E25").Value = "'  +2.40%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'0.1248"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  +4.57%  "
$ws.Range("E26").ClearFormats()
$ws.Range("E27").Value = "'  +1.52%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'0.06469"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  +2.63%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'1.357"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  +4.91%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'1.325"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  +0.69%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'3.603"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  +2.27%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'3.594"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  +2.68%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'1.661"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  +1.67%  "
$ws.Range("E33").ClearFormats()
$ws.Range("E34").Value = "'  +1.96%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'0.6247"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  +3.01%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'2.405"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  +1.77%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'2.719"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  +2.82%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'6.459"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  +4.85%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'1.111.82"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  +3.41%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'0.01625"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  +1.26%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.8779"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  +1.61%  "
$ws.Range("E41").ClearFormats()
$ws.Range("E42").Value = "'  +0.58%  "
$ws.Range("E42").ClearFormats()
$ws.Range("E43").Value = "'  -0.33%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'1.832.01"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  +0.95%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.00000000111"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  -2.09%  "
$ws.Range("E45").ClearFormats()
$ws.Range("E46").Value = "'  +1.75%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'8.210"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  +2.07%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'1.008"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +0.02%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'0.05264"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +1.30%  "
$ws.Range("E49").ClearFormats()
$ws.Range("E50").Value = "'  -0.04%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'6.073"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  +2.70%  "
$ws.Range("E51").ClearFormats()
